# Daily refresh of the cryptos list (prices + 1h volume deltas), plus a
# swap of the dogwifhat/Hedera rows (42/43) to reflect their new rank order.
# Price cells get NumberFormat "@" first so digit strings like "183.80" or
# "0.0000116" are stored verbatim as text instead of being reparsed as
# numbers (which would drop trailing zeros / switch to scientific notation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.085.78"
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.260.07"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.04"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.80"
$ws.Range("E6").Value = "  -1.48%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("E9").Value = "  -3.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.61"
$ws.Range("E11").Value = "  -3.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.825.43"
$ws.Range("E12").Value = "  -0.80%  "
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "68.060.83"
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.24"
$ws.Range("E15").Value = "  -4.17%  "
$ws.Range("E16").Value = "  -2.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.251.40"
$ws.Range("E17").Value = "  -1.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.70"
$ws.Range("E18").Value = "  -3.16%  "
$ws.Range("E19").Value = "  -3.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "417.16"
$ws.Range("E20").Value = "  +5.12%  "
$ws.Range("E21").Value = "  -3.40%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.11"
$ws.Range("E23").Value = "  -0.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.506"
$ws.Range("E24").Value = "  -2.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000116"
$ws.Range("E25").Value = "  -3.84%  "
$ws.Range("E26").Value = "  -1.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.30"
$ws.Range("E27").Value = "  -5.22%  "
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("E29").Value = "  -2.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.57"
$ws.Range("E30").Value = "  -2.48%  "
$ws.Range("E31").Value = "  -5.96%  "
$ws.Range("E32").Value = "  -5.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.23"
$ws.Range("E33").Value = "  -5.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "164.31"
$ws.Range("E34").Value = "  +0.66%  "
$ws.Range("E35").Value = "  -5.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.87"
$ws.Range("E36").Value = "  -7.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.56"
$ws.Range("E37").Value = "  -0.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.791"
$ws.Range("E38").Value = "  -4.25%  "
$ws.Range("E39").Value = "  -4.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.30"
$ws.Range("E40").Value = "  -4.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.621.79"
$ws.Range("E41").Value = "  -1.17%  "
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0670"
$ws.Range("E42").Value = "  -3.38%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.40"
$ws.Range("E43").Value = "  -5.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "335.46"
$ws.Range("E44").Value = "  -2.12%  "
$ws.Range("E45").Value = "  -6.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0272"
$ws.Range("E46").Value = "  -4.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.21"
$ws.Range("E47").Value = "  -2.64%  "
$ws.Range("E48").Value = "  -2.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0999"
$ws.Range("E49").Value = "  -2.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "30.41"
$ws.Range("E51").Value = "  -5.01%  "
